$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: new header cells D1, E1 (same header style as A1:C1) ---
$ws.Cells.Item(1, 4).Value = "PortfolioName"
$ws.Cells.Item(1, 5).Value = "Portfoliocriteria"
$ws.Range("A1").Copy()
$ws.Range("D1:E1").PasteSpecial(-4122)

# --- Row 2: new empty (but present) cells D2, E2 ---
$ws.Cells.Item(2, 4).Font.Bold = $false
$ws.Cells.Item(2, 5).Font.Bold = $false

# --- Row 3: new row, some blank cells, some filled ---
$ws.Cells.Item(3, 1).Font.Bold = $false
$ws.Cells.Item(3, 2).Value = "Arcelo"
$ws.Cells.Item(3, 3).Font.Bold = $false
$ws.Cells.Item(3, 4).Value = "Arcel"
$ws.Cells.Item(3, 5).Value = "Arcelo"
